$wb = $excel.ActiveWorkbook

# Insert the new "2022-Q1" sheet by duplicating the "2021-Q4" sheet's
# layout/formatting, placed immediately before "总计".
$zongjiRef = $wb.Worksheets.Item("总计")
$targetIndex = $zongjiRef.Index
$wb.Worksheets.Item("2021-Q4").Copy($zongjiRef)
$newSheet = $wb.Worksheets.Item($targetIndex)
$newSheet.Name = "2022-Q1"

$data = @(
    @("010861", "长信企业优选一年持有期灵活配置混合", "9.39", "80.21", "3.74", "0.3512", 3),
    @("160211", "国泰中小盘成长混合(LOF)", "6.78", "89.07", "3.67", "0.2488", 6),
    @("005589", "长信企业精选两年定期开放灵活配置混合", "5.84", "79.99", "3.67", "0.2143", 2),
    @("007518", "东方阿尔法优选混合A", "2.03", "72.64", "1.42", "0.0288", 8),
    @("007519", "东方阿尔法优选混合C", "0.82", "72.64", "1.42", "0.0116", 8)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Range("A" + $r).Value = ($r - 2)
    $newSheet.Range("B" + $r).Value = "'" + $row[0]
    $newSheet.Range("C" + $r).Value = $row[1]
    $newSheet.Range("D" + $r).Value = "'" + $row[2]
    $newSheet.Range("E" + $r).Value = "'" + $row[3]
    $newSheet.Range("F" + $r).Value = "'" + $row[4]
    $newSheet.Range("G" + $r).Value = "'" + $row[5]
    $newSheet.Range("H" + $r).Value = $row[6]
    $r = $r + 1
}

# Row 6 is new (source sheet only had 4 data rows) - copy formatting from row 5.
$newSheet.Range("A5").Copy()
$newSheet.Range("A6").PasteSpecial(-4122)

# Update "总计": insert the 2022-Q1 summary row at the top, shifting the
# existing rows (and their running index in column A) down by one.
$zongji = $wb.Worksheets.Item("总计")
for ($r = 7; $r -ge 3; $r--) {
    $zongji.Range("A" + $r).Value = $r - 2
    $zongji.Range("B" + $r).Value = $zongji.Range("B" + ($r - 1)).Value()
    $zongji.Range("C" + $r).Value = $zongji.Range("C" + ($r - 1)).Value()
    $zongji.Range("D" + $r).Value = $zongji.Range("D" + ($r - 1)).Value()
}
$zongji.Range("A6").Copy()
$zongji.Range("A7").PasteSpecial(-4122)

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q1"
$zongji.Range("C2").Value = 5
$zongji.Range("D2").Value = 0.85
